# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.989.86"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.554.66"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "1.776.47"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "1.555.11"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "26.972.64"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "1.405.19"
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.950"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.43%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "1.690.53"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0957"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  +0.52%  "
